$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# November 2023 statistics figures: Circulation, ILL Loans, ILL Borrows per library (rows 3-59)
$data = @(
  @(3, 66465, 10074, 11626),
  @(4, 35529, 4062, 4104),
  @(5, 113522, 10644, 9947),
  @(6, 2542, 1032, 245),
  @(7, 73496, 12208, 9651),
  @(8, 7838, 2026, 1652),
  @(9, 8777, 1711, 963),
  @(10, 4192, 624, 392),
  @(11, 436, 385, 32),
  @(12, 3, 0, 0),
  @(13, 1599, 380, 437),
  @(14, 4745, 1932, 1540),
  @(15, 7539, 2932, 1229),
  @(16, 5532, 2059, 922),
  @(17, 3268, 1355, 272),
  @(18, 26471, 4163, 5051),
  @(19, 2273, 966, 591),
  @(20, 28100, 3576, 4743),
  @(21, 460, 587, 35),
  @(22, 26056, 3622, 4488),
  @(23, 1681, 729, 274),
  @(24, 30362, 3946, 5910),
  @(25, 119562, 11385, 14513),
  @(26, 9282, 3137, 1522),
  @(27, 0, 0, 0),
  @(28, 8141, 1736, 1930),
  @(29, 2170, 674, 456),
  @(30, 22431, 4140, 4018),
  @(31, 679, 257, 325),
  @(32, 4132, 2666, 627),
  @(33, 24709, 4923, 4304),
  @(34, 15704, 4492, 3330),
  @(35, 8640, 976, 1986),
  @(36, 86652, 8568, 8672),
  @(37, 12815, 4190, 1844),
  @(38, 38940, 2960, 4284),
  @(39, 1501, 1435, 251),
  @(40, 2940, 719, 1047),
  @(41, 4510, 549, 174),
  @(42, 16418, 868, 526),
  @(43, 403, 170, 77),
  @(44, 1369, 136, 129),
  @(45, 1045, 14, 7),
  @(46, 5046, 1416, 665),
  @(47, 19067, 5361, 3345),
  @(48, 46739, 5371, 6629),
  @(49, 22345, 5289, 1981),
  @(50, 17168, 1943, 2701),
  @(51, 48453, 4566, 7396),
  @(52, 7308, 1083, 1762),
  @(53, 20639, 4472, 3688),
  @(54, 3067, 2071, 1126),
  @(55, 3506, 1951, 205),
  @(56, 6246, 1639, 2054),
  @(57, 19415, 7303, 4128),
  @(58, 23745, 1737, 913),
  @(59, 1032513, 158170, 145678)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Set zoom level for the active sheet view (85%)
$excel.ActiveWindow.Zoom = 85

Write-Host "Applied November 2023 statistics update"
